$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings remain stored as text (matches source formatting)
$textCells = @("D4", "D5", "D6", "D9", "D10", "D11", "D12", "D14", "D18", "D19", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D31", "D34", "D36", "D40", "D41", "D42", "D45", "D46", "D47", "D48", "D49")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated cell values
$ws.Range('D2').Value = '70.445.94'
$ws.Range('D3').Value = '3.808.75'
$ws.Range('E3').Value = '  +0.73%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '667.94'
$ws.Range('E5').Value = '  +7.26%  '
$ws.Range('D6').Value = '168.62'
$ws.Range('E6').Value = '  +1.77%  '
$ws.Range('D7').Value = '3.806.65'
$ws.Range('E7').Value = '  +0.74%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').Value = '0.527'
$ws.Range('E9').Value = '  +1.09%  '
$ws.Range('D10').Value = '0.161'
$ws.Range('E10').Value = '  +0.05%  '
$ws.Range('D11').Value = '0.463'
$ws.Range('E11').Value = '  +2.44%  '
$ws.Range('D12').Value = '7.04'
$ws.Range('E12').Value = '  +5.92%  '
$ws.Range('E13').Value = '  -1.25%  '
$ws.Range('D14').Value = '35.80'
$ws.Range('E14').Value = '  +0.33%  '
$ws.Range('D15').Value = '4.450.19'
$ws.Range('E15').Value = '  +0.83%  '
$ws.Range('D16').Value = '3.804.86'
$ws.Range('E16').Value = '  +1.03%  '
$ws.Range('D17').Value = '70.371.10'
$ws.Range('E17').Value = '  +1.54%  '
$ws.Range('D18').Value = '17.71'
$ws.Range('D19').Value = '7.16'
$ws.Range('E19').Value = '  +0.77%  '
$ws.Range('E20').Value = '  +0.65%  '
$ws.Range('D21').Value = '10.86'
$ws.Range('E21').Value = '  +12.73%  '
$ws.Range('D22').Value = '474.70'
$ws.Range('E22').Value = '  +1.35%  '
$ws.Range('D23').Value = '0.713'
$ws.Range('E23').Value = '  +1.31%  '
$ws.Range('D24').Value = '82.77'
$ws.Range('E24').Value = '  -0.74%  '
$ws.Range('D25').Value = '0.0000144'
$ws.Range('E25').Value = '  -3.98%  '
$ws.Range('D26').Value = '12.23'
$ws.Range('E26').Value = '  +1.76%  '
$ws.Range('D27').Value = '10.33'
$ws.Range('E27').Value = '  +2.86%  '
$ws.Range('E28').Value = '  -1.87%  '
$ws.Range('E29').Value = '  -0.03%  '
$ws.Range('D30').Value = '3.959.36'
$ws.Range('E30').Value = '  +0.75%  '
$ws.Range('D31').Value = '2.86'
$ws.Range('E31').Value = '  +7.62%  '
$ws.Range('E32').Value = '  +3.13%  '
$ws.Range('E33').Value = '  +0.97%  '
$ws.Range('D34').Value = '29.66'
$ws.Range('E34').Value = '  +2.77%  '
$ws.Range('E35').Value = '  +10.63%  '
$ws.Range('D36').Value = '9.15'
$ws.Range('E36').Value = '  +1.86%  '
$ws.Range('E37').Value = '  +0.02%  '
$ws.Range('D38').Value = '3.766.03'
$ws.Range('E38').Value = '  +0.93%  '
$ws.Range('E39').Value = '  +0.13%  '
$ws.Range('D40').Value = '3.43'
$ws.Range('E40').Value = '  +1.08%  '
$ws.Range('D41').Value = '5.97'
$ws.Range('E41').Value = '  +2.62%  '
$ws.Range('D42').Value = '0.969'
$ws.Range('E42').Value = '  +0.05%  '
$ws.Range('E43').Value = '  +0.07%  '
$ws.Range('D45').Value = '2.09'
$ws.Range('E45').Value = '  +9.49%  '
$ws.Range('D46').Value = '45.84'
$ws.Range('E46').Value = '  +5.97%  '
$ws.Range('B47').Value = 'Monero'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D47').Value = '156.98'
$ws.Range('E47').Value = '  +1.69%  '
$ws.Range('B48').Value = 'OKB'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D48').Value = '48.13'
$ws.Range('E48').Value = '  +2.81%  '
$ws.Range('D49').Value = '0.301'
$ws.Range('E49').Value = '  +0.35%  '
$ws.Range('E50').Value = '  +4.14%  '
$ws.Range('E51').Value = '  +0.85%  '
